$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.639.72"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "3.445.05"
$ws.Range("E3").Value = "  -1.20%  "

$ws.Range("E4").Value = "  -0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.06"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "175.01"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.38%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.598"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.43%  "

$ws.Range("D9").Value = "3.444.84"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("E10").Value = "  -3.14%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.83"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.61%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.420"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.60%  "

$ws.Range("D13").Value = "4.047.08"
$ws.Range("E13").Value = "  -1.09%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "30.94"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.05%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.131"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.49%  "

$ws.Range("D16").Value = "66.623.73"
$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("E17").Value = "  -3.21%  "

$ws.Range("D18").Value = "3.442.75"
$ws.Range("E18").Value = "  -1.41%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.02"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.80%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.79"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.62%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "375.36"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.04%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.70"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.26%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.25%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "70.84"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.95%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.525"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.96%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0000118"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.87"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.95%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.172"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "

$ws.Range("E30").Value = "  -0.15%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.86"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.59%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "23.88"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.27%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.11%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.33"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.51%  "

$ws.Range("E35").Value = "  -0.08%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "7.05"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.19%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.51"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.58%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "159.01"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.14%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.877"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "27.12"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.96%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.78"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.79%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.61"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.23%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.49"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -5.37%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.43"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.18%  "

$ws.Range("D45").Value = "2.686.69"
$ws.Range("E45").Value = "  -5.58%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0690"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.64%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "25.17"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.01%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "40.40"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.86%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0293"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.97%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "319.13"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -5.46%  "

$ws.Range("E51").Value = "  -4.45%  "
